$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text type for numeric-looking strings in column D by temporarily
# forcing Text format, then clearing the format back afterwards so the
# cells retain their original (unstyled) appearance while keeping t="s".
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = "42.354.69"
$ws.Range("E2").Value2 = "  +0.35%  "
$ws.Range("D3").Value2 = "2.299.79"
$ws.Range("E3").Value2 = "  -0.64%  "
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").Value2 = "316.55"
$ws.Range("E5").Value2 = "  +1.24%  "
$ws.Range("D6").Value2 = "103.40"
$ws.Range("E6").Value2 = "  -2.75%  "
$ws.Range("E7").Value2 = "  +0.56%  "
$ws.Range("E8").Value2 = "  -0.10%  "
$ws.Range("D9").Value2 = "0.610"
$ws.Range("E9").Value2 = "  -0.13%  "
$ws.Range("D10").Value2 = "39.94"
$ws.Range("E10").Value2 = "  -0.74%  "
$ws.Range("D11").Value2 = "0.0910"
$ws.Range("E11").Value2 = "  -0.78%  "
$ws.Range("D12").Value2 = "8.32"
$ws.Range("E12").Value2 = "  +0.16%  "
$ws.Range("E13").Value2 = "  +0.09%  "
$ws.Range("E14").Value2 = "  -1.47%  "
$ws.Range("D15").Value2 = "15.25"
$ws.Range("E15").Value2 = "  -2.08%  "
$ws.Range("D16").Value2 = "2.648.11"
$ws.Range("E16").Value2 = "  -0.44%  "
$ws.Range("D17").Value2 = "2.287.57"
$ws.Range("E17").Value2 = "  -1.44%  "
$ws.Range("D18").Value2 = "42.474.12"
$ws.Range("E18").Value2 = "  +0.74%  "
$ws.Range("D19").Value2 = "7.49"
$ws.Range("E19").Value2 = "  -3.38%  "
$ws.Range("E20").Value2 = "  +0.69%  "
$ws.Range("D21").Value2 = "72.99"
$ws.Range("E21").Value2 = "  -2.27%  "
$ws.Range("B22").Value2 = "BitcoinCash"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value2 = "278.42"
$ws.Range("E22").Value2 = "  +7.37%  "
$ws.Range("B23").Value2 = "PancakeSwap"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value2 = "3.56"
$ws.Range("E23").Value2 = "  +2.11%  "
$ws.Range("D24").Value2 = "11.12"
$ws.Range("E24").Value2 = "  +19.92%  "
$ws.Range("D25").Value2 = "2.27"
$ws.Range("E25").Value2 = "  -0.85%  "
$ws.Range("E26").Value2 = "  -0.35%  "
$ws.Range("D27").Value2 = "10.84"
$ws.Range("E27").Value2 = "  -1.87%  "
$ws.Range("E28").Value2 = "  +3.27%  "
$ws.Range("D29").Value2 = "22.75"
$ws.Range("E29").Value2 = "  -0.34%  "
$ws.Range("D30").Value2 = "35.81"
$ws.Range("E30").Value2 = "  +0.61%  "
$ws.Range("D31").Value2 = "165.51"
$ws.Range("E31").Value2 = "  +1.73%  "
$ws.Range("D32").Value2 = "0.0874"
$ws.Range("E32").Value2 = "  -2.25%  "
$ws.Range("D33").Value2 = "5.87"
$ws.Range("E33").Value2 = "  -0.12%  "
$ws.Range("E34").Value2 = "  +5.42%  "
$ws.Range("D35").Value2 = "0.118"
$ws.Range("E35").Value2 = "  +0.34%  "
$ws.Range("D36").Value2 = "2.60"
$ws.Range("E36").Value2 = "  -10.75%  "
$ws.Range("D37").Value2 = "0.0370"
$ws.Range("E37").Value2 = "  +4.47%  "
$ws.Range("E38").Value2 = "  +1.86%  "
$ws.Range("D39").Value2 = "3.73"
$ws.Range("E39").Value2 = "  +2.51%  "
$ws.Range("E40").Value2 = "  +0.53%  "
$ws.Range("E41").Value2 = "  +1.87%  "
$ws.Range("D42").Value2 = "96.19"
$ws.Range("E42").Value2 = "  -2.01%  "
$ws.Range("D43").Value2 = "69.72"
$ws.Range("E43").Value2 = "  -1.40%  "
$ws.Range("D44").Value2 = "0.227"
$ws.Range("E44").Value2 = "  -2.07%  "
$ws.Range("E45").Value2 = "  +0.26%  "
$ws.Range("D46").Value2 = "82.24"
$ws.Range("E46").Value2 = "  +10.25%  "
$ws.Range("D47").Value2 = "12.05"
$ws.Range("E47").Value2 = "  -0.92%  "
$ws.Range("D48").Value2 = "112.24"
$ws.Range("E48").Value2 = "  +0.31%  "
$ws.Range("E49").Value2 = "  -0.28%  "
$ws.Range("B50").Value2 = "Maker"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value2 = "1.592.01"
$ws.Range("E50").Value2 = "  +3.14%  "
$ws.Range("B51").Value2 = "THORChain"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value2 = "5.16"
$ws.Range("E51").Value2 = "  -4.29%  "

$ws.Range("D2:D51").ClearFormats()
